$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "It should not be displayed and should be as per parameter."

# Rows 19-21 (TC02 vendor/payableID/entryID/vendorID block) and
# rows 46-48 (TC05 equivalent block): text updated, same row height.
$ws.Range("F19").Value = $newText
$ws.Range("F20").Value = $newText
$ws.Range("F21").Value = $newText

$ws.Range("F46").Value = $newText
$ws.Range("F47").Value = $newText
$ws.Range("F48").Value = $newText

# Rows 68-70 (TC07 block): text updated AND row height reverts to the
# sheet's default (the explicit 29pt height is cleared).
$ws.Range("F68").Value = $newText
$ws.Range("F69").Value = $newText
$ws.Range("F70").Value = $newText
$ws.Range("A68:F70").Rows.AutoFit()

# Match the saved selection/active cell left in the sheet view.
$ws.Range("F68:F70").Select()
